$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Replace 26x80= -> 98x15= FIRST, before 21x90= -> 26x80= is created,
# to avoid cascading replacement of the newly-created 26x80=.
Replace-Text "26×80=" "98×15="

Replace-Text "48×38=" "96×71="
Replace-Text "54×19=" "97×93="
Replace-Text "45×99=" "98×48="
Replace-Text "12×57=" "57×53="
Replace-Text "54×73=" "44×15="
Replace-Text "62×80=" "54×69="
Replace-Text "67×31=" "45×34="
Replace-Text "84×47=" "17×93="
Replace-Text "71×24=" "45×30="
Replace-Text "21×90=" "26×80="
Replace-Text "43×72=" "61×91="
Replace-Text "55×62=" "88×98="
Replace-Text "70×40=" "77×99="
Replace-Text "50×72=" "33×59="
Replace-Text "34×71=" "85×68="
Replace-Text "16×71=" "80×16="
Replace-Text "67×95=" "57×24="
Replace-Text "50×71=" "42×81="
Replace-Text "16×92=" "69×80="
Replace-Text "42×52=" "78×12="
Replace-Text "33×27=" "25×21="
Replace-Text "19×26=" "70×55="
Replace-Text "21×88=" "78×44="
Replace-Text "41×36=" "92×74="
